# Update crypto price/volume table to refreshed values (GitHub Actions data pull).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = '62.146.39'
$cell.ClearFormats()
$ws.Cells.Item(2, 5).Value = '  +1.45%  '

$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.424.98'
$cell.ClearFormats()
$ws.Cells.Item(3, 5).Value = '  +0.77%  '

$ws.Cells.Item(4, 5).Value = '  -0.18%  '

$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '406.79'
$cell.ClearFormats()
$ws.Cells.Item(5, 5).Value = '  +0.26%  '

$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = '131.20'
$cell.ClearFormats()
$ws.Cells.Item(6, 5).Value = '  +1.98%  '

$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.595'
$cell.ClearFormats()
$ws.Cells.Item(7, 5).Value = '  -2.06%  '

$ws.Cells.Item(8, 5).Value = '  -0.17%  '

$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.691'
$cell.ClearFormats()
$ws.Cells.Item(9, 5).Value = '  +2.81%  '

$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.138'
$cell.ClearFormats()
$ws.Cells.Item(10, 5).Value = '  +9.17%  '

$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = '41.80'
$cell.ClearFormats()
$ws.Cells.Item(11, 5).Value = '  -1.10%  '

$ws.Cells.Item(12, 5).Value = '  -0.22%  '

$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = '19.89'
$cell.ClearFormats()
$ws.Cells.Item(13, 5).Value = '  +1.18%  '

$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = '8.41'
$cell.ClearFormats()
$ws.Cells.Item(14, 5).Value = '  -1.60%  '

$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.440.35'
$cell.ClearFormats()
$ws.Cells.Item(15, 5).Value = '  +1.04%  '

$ws.Cells.Item(16, 2).Value = 'Uniswap'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = '11.59'
$cell.ClearFormats()
$ws.Cells.Item(16, 5).Value = '  -0.46%  '

$ws.Cells.Item(17, 2).Value = 'WrappedBTC'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = '62.023.63'
$cell.ClearFormats()
$ws.Cells.Item(17, 5).Value = '  +1.05%  '

$ws.Cells.Item(18, 5).Value = '  -0.13%  '

$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0000150'
$cell.ClearFormats()
$ws.Cells.Item(19, 5).Value = '  +11.15%  '

$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.17'
$cell.ClearFormats()
$ws.Cells.Item(20, 5).Value = '  -2.39%  '

$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = '84.17'
$cell.ClearFormats()
$ws.Cells.Item(21, 5).Value = '  +1.77%  '

$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = '312.26'
$cell.ClearFormats()
$ws.Cells.Item(22, 5).Value = '  +1.59%  '

$ws.Cells.Item(23, 5).Value = '  -2.69%  '

$ws.Cells.Item(24, 5).Value = '  +0.20%  '

$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.75'
$cell.ClearFormats()
$ws.Cells.Item(25, 5).Value = '  +0.63%  '

$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = '29.69'
$cell.ClearFormats()
$ws.Cells.Item(26, 5).Value = '  +0.15%  '

$ws.Cells.Item(27, 2).Value = 'Filecoin'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = '8.17'
$cell.ClearFormats()
$ws.Cells.Item(27, 5).Value = '  -4.95%  '

$ws.Cells.Item(28, 2).Value = 'RenderToken'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.84'
$cell.ClearFormats()
$ws.Cells.Item(28, 5).Value = '  +4.37%  '

$ws.Cells.Item(29, 5).Value = '  +5.35%  '

$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.172'
$cell.ClearFormats()
$ws.Cells.Item(30, 5).Value = '  -0.97%  '

$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = '43.68'
$cell.ClearFormats()
$ws.Cells.Item(31, 5).Value = '  +1.61%  '

$ws.Cells.Item(32, 5).Value = '  -1.04%  '

$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = '11.32'
$cell.ClearFormats()
$ws.Cells.Item(33, 5).Value = '  -3.33%  '

$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.ClearFormats()
$ws.Cells.Item(34, 5).Value = '  -0.08%  '

$ws.Cells.Item(35, 5).Value = '  -0.29%  '

$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = '51.58'
$cell.ClearFormats()
$ws.Cells.Item(36, 5).Value = '  -0.99%  '

$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.997'
$cell.ClearFormats()
$ws.Cells.Item(37, 5).Value = '  -0.08%  '

$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.02'
$cell.ClearFormats()
$ws.Cells.Item(38, 5).Value = '  +0.92%  '

$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.33'
$cell.ClearFormats()
$ws.Cells.Item(39, 5).Value = '  -2.97%  '

$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.314'
$cell.ClearFormats()
$ws.Cells.Item(40, 5).Value = '  +10.04%  '

$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = '140.38'
$cell.ClearFormats()
$ws.Cells.Item(41, 5).Value = '  +2.89%  '

$ws.Cells.Item(42, 5).Value = '  -0.23%  '

$ws.Cells.Item(43, 5).Value = '  +0.04%  '

$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.93'
$cell.ClearFormats()
$ws.Cells.Item(44, 5).Value = '  -0.25%  '

$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = '16.75'
$cell.ClearFormats()
$ws.Cells.Item(45, 5).Value = '  -1.21%  '

$ws.Cells.Item(46, 5).Value = '  -0.05%  '

$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = '21.31'
$cell.ClearFormats()
$ws.Cells.Item(47, 5).Value = '  -2.43%  '

$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.099.61'
$cell.ClearFormats()
$ws.Cells.Item(48, 5).Value = '  -2.24%  '

$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.32'
$cell.ClearFormats()
$ws.Cells.Item(49, 5).Value = '  -0.44%  '

$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.95'
$cell.ClearFormats()
$ws.Cells.Item(50, 5).Value = '  +1.57%  '

$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.70'
$cell.ClearFormats()
$ws.Cells.Item(51, 5).Value = '  +15.81%  '

